$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 36
$ws.Range("H2").Value = 48
$ws.Range("E4").Value = 26
$ws.Range("E7").Value = 10
$ws.Range("E25").Value = 23
$ws.Range("E26").Value = 32
$ws.Range("E34").Value = 25
$ws.Range("E35").Value = 11
$ws.Range("E38").Value = 81
$ws.Range("E43").Value = 27
$ws.Range("E45").Value = 27
$ws.Range("E56").Value = 9
$ws.Range("E57").Value = 16
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = 5
$ws.Range("H59").Value = 9
$ws.Range("E63").Value = 41
$ws.Range("F63").Value = 14
$ws.Range("H63").Value = 22
$ws.Range("E71").Value = 43
$ws.Range("F71").Value = 21
$ws.Range("H71").Value = 31
$ws.Range("E73").Value = 31
$ws.Range("E77").Value = 58
$ws.Range("F77").Value = 23
$ws.Range("H77").Value = 40
$ws.Range("E80").Value = 31
$ws.Range("F80").Value = 15
$ws.Range("H80").Value = 27
$ws.Range("E87").Value = 19
